# Scheduled data refresh: update currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) across the per-job sheets with latest market-board pulls.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1603.8889
$ws.Range("I40").Value = 1610.1666
$ws.Range("J40").Value = 1591.3334
$ws.Range("K40").Value = 1610.1666
$ws.Range("L40").Value = 1591.3334
$ws.Range("M40").Value = -1435.1666
$ws.Range("N40").Value = -1941.3334

$ws.Range("H64").Value = 396225.38
$ws.Range("I64").Value = 641203.8
$ws.Range("K64").Value = 641203.8
$ws.Range("M64").Value = -640955.8

$ws.Range("H67").Value = 396225.38
$ws.Range("I67").Value = 641203.8
$ws.Range("K67").Value = 641203.8
$ws.Range("M67").Value = -640345.8

$ws.Range("H101").Value = 1530.6666
$ws.Range("I101").Value = 362.36365
$ws.Range("J101").Value = 2815.8
$ws.Range("K101").Value = 1087.09095
$ws.Range("L101").Value = 8447.400000000001
$ws.Range("M101").Value = 534.90905
$ws.Range("N101").Value = -11691.4

$ws.Range("H107").Value = 1012.5
$ws.Range("I107").Value = 866.25
$ws.Range("K107").Value = 866.25
$ws.Range("M107").Value = 1053.75

$ws.Range("H132").Value = 2081.3547
$ws.Range("I132").Value = 2052.4827
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 6157.4481
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3627.4481
$ws.Range("N132").Value = -12560

$ws.Range("H138").Value = 2062.9495
$ws.Range("I138").Value = 1650.875
$ws.Range("J138").Value = 2142.3855
$ws.Range("K138").Value = 4952.625
$ws.Range("L138").Value = 6427.156499999999
$ws.Range("M138").Value = 187.375
$ws.Range("N138").Value = -16707.1565

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1040.9678
$ws.Range("I2").Value = 1085.1666
$ws.Range("J2").Value = 889.4286
$ws.Range("K2").Value = 1085.1666
$ws.Range("L2").Value = 889.4286
$ws.Range("M2").Value = -972.1666
$ws.Range("N2").Value = -1115.4286

$ws.Range("H45").Value = 2062.087
$ws.Range("I45").Value = 1990.2106
$ws.Range("J45").Value = 2403.5
$ws.Range("K45").Value = 1990.2106
$ws.Range("L45").Value = 2403.5
$ws.Range("M45").Value = -1613.2106
$ws.Range("N45").Value = -3157.5

$ws.Range("H63").Value = 3399.6924
$ws.Range("I63").Value = 3472.3635
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3472.3635
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2786.3635
$ws.Range("N63").Value = -4372

$ws.Range("H66").Value = 3399.6924
$ws.Range("I66").Value = 3472.3635
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 17361.8175
$ws.Range("L66").Value = 3000
$ws.Range("M66").Value = -13929.8175
$ws.Range("N66").Value = -21864

$ws.Range("H74").Value = 7279.5
$ws.Range("I74").Value = 7279.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 7279.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -6405.5
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 7279.5
$ws.Range("I77").Value = 7279.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 36397.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -32029.5
$ws.Range("N77").ClearContents()

$ws.Range("H116").Value = 1040.9678
$ws.Range("I116").Value = 1085.1666
$ws.Range("J116").Value = 889.4286
$ws.Range("K116").Value = 1085.1666
$ws.Range("L116").Value = 889.4286
$ws.Range("M116").Value = 1208.8334
$ws.Range("N116").Value = -5477.4286

$ws.Range("H122").Value = 1512.5483
$ws.Range("I122").Value = 1486.4
$ws.Range("J122").Value = 1621.5
$ws.Range("K122").Value = 4459.200000000001
$ws.Range("L122").Value = 4864.5
$ws.Range("M122").Value = -2009.200000000001
$ws.Range("N122").Value = -9764.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1040.9678
$ws.Range("I3").Value = 1085.1666
$ws.Range("J3").Value = 889.4286
$ws.Range("K3").Value = 1085.1666
$ws.Range("L3").Value = 889.4286
$ws.Range("M3").Value = -971.1666
$ws.Range("N3").Value = -1117.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 259.33334
$ws.Range("I22").Value = 80
$ws.Range("K22").Value = 80
$ws.Range("M22").Value = 270

$ws.Range("H57").Value = 21212
$ws.Range("J57").Value = 21212
$ws.Range("L57").Value = 21212
$ws.Range("N57").Value = -22332

$ws.Range("H122").Value = 8000
$ws.Range("I122").Value = 9500
$ws.Range("K122").Value = 28500
$ws.Range("M122").Value = -26050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 70000
$ws.Range("J37").Value = 70000
$ws.Range("L37").Value = 210000
$ws.Range("N37").Value = -210224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2882.5217
$ws.Range("I126").Value = 1937.25
$ws.Range("J126").Value = 3386.6667
$ws.Range("K126").Value = 5811.75
$ws.Range("L126").Value = 10160.0001
$ws.Range("M126").Value = -3341.75
$ws.Range("N126").Value = -15100.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5533.1904
$ws.Range("I7").Value = 4086.2
$ws.Range("K7").Value = 4086.2
$ws.Range("M7").Value = -3974.2

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H61").Value = 11532.444
$ws.Range("I61").Value = 14000.363
$ws.Range("J61").Value = 7654.2856
$ws.Range("K61").Value = 14000.363
$ws.Range("L61").Value = 7654.2856
$ws.Range("M61").Value = -13798.363
$ws.Range("N61").Value = -8058.2856

$ws.Range("H68").Value = 1834.1666
$ws.Range("I68").Value = 1375.5
$ws.Range("J68").Value = 2751.5
$ws.Range("K68").Value = 1375.5
$ws.Range("L68").Value = 2751.5
$ws.Range("M68").Value = -626.5
$ws.Range("N68").Value = -4249.5

$ws.Range("H71").Value = 1834.1666
$ws.Range("I71").Value = 1375.5
$ws.Range("J71").Value = 2751.5
$ws.Range("K71").Value = 6877.5
$ws.Range("L71").Value = 13757.5
$ws.Range("M71").Value = -3133.5
$ws.Range("N71").Value = -21245.5

$ws.Range("H107").Value = 3466.6667
$ws.Range("I107").Value = 3466.6667
$ws.Range("K107").Value = 3466.6667
$ws.Range("M107").Value = -1546.6667

$ws.Range("H113").Value = 11532.444
$ws.Range("I113").Value = 14000.363
$ws.Range("J113").Value = 7654.2856
$ws.Range("K113").Value = 14000.363
$ws.Range("L113").Value = 7654.2856
$ws.Range("M113").Value = -11830.363
$ws.Range("N113").Value = -11994.2856

$ws.Range("H126").Value = 5533.1904
$ws.Range("I126").Value = 4086.2
$ws.Range("K126").Value = 12258.6
$ws.Range("M126").Value = -9788.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1510
$ws.Range("I107").Value = 881.7857
$ws.Range("J107").Value = 2389.5
$ws.Range("K107").Value = 2645.3571
$ws.Range("L107").Value = 7168.5
$ws.Range("M107").Value = -725.3571000000002
$ws.Range("N107").Value = -11008.5

$ws.Range("H113").Value = 1111.5555
$ws.Range("I113").Value = 550.625
$ws.Range("J113").Value = 1560.3
$ws.Range("K113").Value = 1651.875
$ws.Range("L113").Value = 4680.9
$ws.Range("M113").Value = 518.125
$ws.Range("N113").Value = -9020.9

$ws.Range("H122").Value = 1996.5897
$ws.Range("I122").Value = 1770.0714
$ws.Range("J122").Value = 2573.182
$ws.Range("K122").Value = 5310.2142
$ws.Range("L122").Value = 7719.545999999999
$ws.Range("M122").Value = -2860.2142
$ws.Range("N122").Value = -12619.546

$ws.Range("H132").Value = 2033.8276
$ws.Range("I132").Value = 969.7692
$ws.Range("J132").Value = 2898.375
$ws.Range("K132").Value = 2909.3076
$ws.Range("L132").Value = 8695.125
$ws.Range("M132").Value = -379.3076000000001
$ws.Range("N132").Value = -13755.125
